# Updated cryptos list (GitHub Actions refresh): new Price (D) / Volume(1h) (E)
# values for each coin row. D-column values are forced to text (leading
# apostrophe) because several look like numbers (e.g. "218.68", "1.012")
# and Excel would otherwise silently convert them to numeric cells; the
# style is then reset to "Normal" so the quote-prefix formatting doesn't
# stick to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.406.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "'1.691.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.90%  "
$ws.Range("D5").Value = "'218.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'0.5471"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.13%  "
$ws.Range("E7").Value = "  +0.84%  "
$ws.Range("D8").Value = "'0.2715"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'0.06464"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "'0.07705"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").Value = "'1.702.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "'4.535"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("D14").Value = "'0.5816"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "'0.000008383"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "'65.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "'26.461.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("D20").Value = "'10.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").Value = "'189.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'1.012"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +3.84%  "
$ws.Range("D25").Value = "'0.1301"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.43%  "
$ws.Range("D26").Value = "'7.877"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.65%  "
$ws.Range("D27").Value = "'15.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'1.420"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.03%  "
$ws.Range("D29").Value = "'0.06320"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.21%  "
$ws.Range("D30").Value = "'1.330"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "'3.578"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "'3.573"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "'1.041"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("D35").Value = "'0.6214"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'2.416"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").Value = "'2.725"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Value = "'6.210"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'1.118.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("D41").Value = "'0.8784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'101.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'1.844.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  -4.94%  "
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "'8.206"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").Value = "'1.008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "'0.05278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Value = "'0.4308"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'6.063"
$ws.Range("D51").Style = "Normal"
